$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Footer "date updates automatically" field: 04.08.2025 -> 06.08.2025
#    on the Slide Master and every Slide Layout (12 occurrences total).
# ---------------------------------------------------------------------
$newDate = "06.08.2025"

$master = $p.SlideMaster
for ($i = 1; $i -le $master.Shapes.Count; $i++) {
    $mshp = $master.Shapes.Item($i)
    if ($mshp.HasTextFrame -and $mshp.Name -like "Date Placeholder*") {
        $mshp.TextFrame.TextRange.Text = $newDate
    }
}

for ($L = 1; $L -le $master.CustomLayouts.Count; $L++) {
    $layout = $master.CustomLayouts.Item($L)
    for ($i = 1; $i -le $layout.Shapes.Count; $i++) {
        $lshp = $layout.Shapes.Item($i)
        if ($lshp.HasTextFrame -and $lshp.Name -like "Date Placeholder*") {
            $lshp.TextFrame.TextRange.Text = $newDate
        }
    }
}

# ---------------------------------------------------------------------
# 2) Slide 2 ("WBZ451H side modifications") picture + two callouts.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)

# Background picture nudged right/down a bit.
$pic = $s2.Shapes.Item(1)
$pic.Left = 669472 / 12700
$pic.Top = 61912 / 12700

# Small highlight rectangle around the UART pins grows/moves.
$rect = $s2.Shapes.Item(2)
$rect.Left = 6801852 / 12700
$rect.Top = 3512676 / 12700
$rect.Width = 535119 / 12700
$rect.Height = 553998 / 12700

# "Remove R4 / to disable RTS..." callout: widened, and its text is
# extended in place (keeping the manual line break between the runs).
$callout = $s2.Shapes.Item(8)
$callout.Left = 6524841 / 12700
$callout.Width = 1889239 / 12700

$tr = $callout.TextFrame.TextRange
$firstRun = $tr.Characters(1, 10)
$firstRun.Text = "Remove R1, R2, R4 "

$tr = $callout.TextFrame.TextRange
$secondRun = $tr.Characters(20, 42)
$secondRun.Text = "to disable RX, TX, RTS from UART Virtual Comm Port"

# ---------------------------------------------------------------------
# 3) Slide 4 pin table: widen it and update two "not connected" cells.
# ---------------------------------------------------------------------
$s4 = $p.Slides.Item(4)
$tblShape = $s4.Shapes.Item(1)
$tblShape.Width = 6467914 / 12700

$tbl = $tblShape.Table
$dash = [char]0x2013
$tbl.Cell(14, 4).Shape.TextFrame.TextRange.Text = "not connected " + $dash + " used as console RX"
$tbl.Cell(15, 4).Shape.TextFrame.TextRange.Text = "not connected " + $dash + " used as console TX"
